# Refresh the cryptos list snapshot (prices / 1h volume deltas) for rows 2-51,
# plus two coin rows (27/28 and 49/50) whose ranking swapped position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "54.714.76"
$ws.Range("E2").Value = "  +0.98%  "

# Row 3
$ws.Range("D3").Value = "2.299.10"
$ws.Range("E3").Value = "  +0.48%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.27%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "498.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.35%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.03%  "

# Row 9
$ws.Range("D9").Value = "2.298.44"
$ws.Range("E9").Value = "  +0.65%  "

# Row 10
$ws.Range("E10").Value = "  +1.87%  "

# Row 11
$ws.Range("E11").Value = "  +2.27%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.325"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.62%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.93%  "

# Row 14
$ws.Range("D14").Value = "2.703.46"
$ws.Range("E14").Value = "  +0.08%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.73%  "

# Row 16
$ws.Range("D16").Value = "54.625.49"
$ws.Range("E16").Value = "  +0.87%  "

# Row 17
$ws.Range("E17").Value = "  +0.81%  "

# Row 18
$ws.Range("D18").Value = "2.337.60"
$ws.Range("E18").Value = "  +0.98%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.77%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "306.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.42%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.44%  "

# Row 24
$ws.Range("E24").Value = "  -1.42%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.31%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.25%  "

# Row 27
$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.375"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.00%  "

# Row 28
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.152"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.27%  "

# Row 29
$ws.Range("D29").Value = "2.396.92"
$ws.Range("E29").Value = "  -1.38%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.87%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.33%  "

# Row 32
$ws.Range("E32").Value = "  +0.08%  "

# Row 33
$ws.Range("D33").Value = "0.0₃0693"
$ws.Range("E33").Value = "  -0.56%  "

# Row 34
$ws.Range("E34").Value = "  +2.85%  "

# Row 35
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("E36").Value = "  +2.12%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.05%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.57%  "

# Row 39
$ws.Range("E39").Value = "  +3.70%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.868"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.34%  "

# Row 41
$ws.Range("E41").Value = "  +1.81%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.75%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.60%  "

# Row 44
$ws.Range("E44").Value = "  +2.65%  "

# Row 45
$ws.Range("E45").Value = "  +1.28%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "129.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.17%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.51%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0895"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.41%  "

# Row 49
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "245.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.09%  "

# Row 50
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.550"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.50%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0486"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.62%  "
